# Generate Report for Handoff
# Updates the localization-status report: refresh "handed back" rows to
# "Ready for handoff", bump the xliff-generation / handoff timestamps,
# flip Priority ht -> mt, and flag the stale fbccde18 handback file with
# an Error Detail message. Also re-autofits a couple of report columns.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$readyForHandoff = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8c0d9baab768449e650b38611646af84d6a5c3a2/e2e/fbccde18-7968-418e-9bfb-cf6fbbfdb8e3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0d8ecff411167679723c9531ac25feca90d99f7/e2e/fbccde18-7968-418e-9bfb-cf6fbbfdb8e3.md."

# ---- Overview sheet ----
$wsOverview.Range("E2").Value = $readyForHandoff
$wsOverview.Range("F2").Value = $readyForHandoff
$wsOverview.Range("E3").Value = $readyForHandoff
$wsOverview.Range("F3").Value = $readyForHandoff
$wsOverview.Range("G2").Value = "2016-12-16 09:56:15"
$wsOverview.Range("G3").Value = "2016-12-16 09:56:15"

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# ---- zh-cn sheet ----
$wsZhCn.Range("C2").Value = $readyForHandoff
$wsZhCn.Range("C3").Value = $readyForHandoff
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-12-16 09:56:01"
$wsZhCn.Range("H3").Value = "2016-12-16 09:56:01"
$wsZhCn.Range("R3").Value = $errorDetail

$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(18).ColumnWidth = 40

# ---- de-de sheet ----
$wsDeDe.Range("C2").Value = $readyForHandoff
$wsDeDe.Range("C3").Value = $readyForHandoff
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-12-16 09:56:15"
$wsDeDe.Range("H3").Value = "2016-12-16 09:56:15"
$wsDeDe.Range("R3").Value = $errorDetail

$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(18).ColumnWidth = 40
